$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Action1")

# A1 held "a"; grow it to "aaaaa" (the only real content change in the diff).
# (Re-setting it as a shared string naturally reshuffles the sharedStrings
# table the same way Excel does: the now-unused "a" entry drops out and the
# new "aaaaa" entry is appended.)
$ws.Range("A1").Value = "aaaaa"

# Move the active selection back to A1 so the sheet view reverts to the
# (unselected / default) state instead of keeping the old B2 selection.
$ws.Activate()
$ws.Range("A1").Select()
